$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F6").Value = 3808
$ws.Range("F8").Value = 2535
$ws.Range("F9").Value = 71
$ws.Range("F10").Value = 3076
$ws.Range("F12").Value = 531
$ws.Range("F13").Value = 2290
$ws.Range("F14").Value = 49
$ws.Range("F16").Value = 86
$ws.Range("F17").Value = 442
$ws.Range("F18").Value = 3
$ws.Range("F20").Value = 201
$ws.Range("F21").Value = 342
$ws.Range("F23").Value = 351
$ws.Range("F24").Value = 647
$ws.Range("F25").Value = 1398
$ws.Range("F26").Value = 42
$ws.Range("F28").Value = 124
$ws.Range("F29").Value = 146
$ws.Range("F30").Value = 1
$ws.Range("F32").Value = 39
$ws.Range("F33").Value = 4240
$ws.Range("F34").Value = 3915
$ws.Range("F35").Value = 70
$ws.Range("F38").Value = 1115
$ws.Range("F40").Value = 464
$ws.Range("F43").Value = 158
$ws.Range("F46").Value = 36
$ws.Range("F48").Value = 54

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 143
$ws.Range("F4").Value = 2267

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 143
$ws.Range("F10").Value = 3808
$ws.Range("F12").Value = 2535
$ws.Range("F13").Value = 71
$ws.Range("F14").Value = 3076
$ws.Range("F15").Value = 531
$ws.Range("F16").Value = 2290
$ws.Range("F17").Value = 49
$ws.Range("F19").Value = 86
$ws.Range("F20").Value = 442
$ws.Range("F22").Value = 342
$ws.Range("F24").Value = 351
$ws.Range("F25").Value = 647
$ws.Range("F26").Value = 1398
$ws.Range("F27").Value = 42
$ws.Range("F29").Value = 124
$ws.Range("F30").Value = 146
$ws.Range("F32").Value = 39
$ws.Range("F35").Value = 4240
$ws.Range("F36").Value = 3915
$ws.Range("F37").Value = 70
$ws.Range("F38").Value = 1115
$ws.Range("F40").Value = 464
$ws.Range("F45").Value = 158
$ws.Range("F48").Value = 54
